# Applies the crypto price/volume refresh captured in the commit:
# "Updated cryptos list on Mon Nov 18 21:59:11 UTC 2024 with GitHub Actions"
#
# Column D ("Price") and E ("Volume(1h)") are plain text cells (values such
# as "91.406.06" or "  +2.23%  " are not real numbers), and a couple of rows
# also had their Coin/Link/Price/Volume swapped with a neighboring row.
# Every update below is expressed as a literal string write.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "91.331.62" }
    @{ Cell = "E2"; Value = "  +2.29%  " }
    @{ Cell = "D3"; Value = "3.150.76" }
    @{ Cell = "E3"; Value = "  +2.71%  " }
    @{ Cell = "E4"; Value = "  +0.03%  " }
    @{ Cell = "D5"; Value = "238.66" }
    @{ Cell = "E5"; Value = "  +2.01%  " }
    @{ Cell = "D6"; Value = "619.27" }
    @{ Cell = "E6"; Value = "  +0.64%  " }
    @{ Cell = "E7"; Value = "  +6.17%  " }
    @{ Cell = "D8"; Value = "0.373" }
    @{ Cell = "E8"; Value = "  +4.34%  " }
    @{ Cell = "D9"; Value = "0.999" }
    @{ Cell = "E9"; Value = "  -0.22%  " }
    @{ Cell = "D10"; Value = "3.149.98" }
    @{ Cell = "E10"; Value = "  +2.75%  " }
    @{ Cell = "D11"; Value = "0.742" }
    @{ Cell = "E11"; Value = "  +5.16%  " }
    @{ Cell = "D12"; Value = "0.202" }
    @{ Cell = "E12"; Value = "  +2.55%  " }
    @{ Cell = "E13"; Value = "  -0.73%  " }
    @{ Cell = "D14"; Value = "35.22" }
    @{ Cell = "E14"; Value = "  +0.99%  " }
    @{ Cell = "D15"; Value = "5.58" }
    @{ Cell = "E15"; Value = "  +4.62%  " }
    @{ Cell = "D16"; Value = "91.461.33" }
    @{ Cell = "E16"; Value = "  +2.42%  " }
    @{ Cell = "D17"; Value = "3.739.75" }
    @{ Cell = "E17"; Value = "  +2.81%  " }
    @{ Cell = "D18"; Value = "3.185.84" }
    @{ Cell = "E18"; Value = "  +4.30%  " }
    @{ Cell = "D19"; Value = "3.73" }
    @{ Cell = "E19"; Value = "  -0.36%  " }
    @{ Cell = "D20"; Value = "15.28" }
    @{ Cell = "E20"; Value = "  +11.77%  " }
    @{ Cell = "D21"; Value = "5.95" }
    @{ Cell = "E21"; Value = "  +11.12%  " }
    @{ Cell = "D22"; Value = "455.10" }
    @{ Cell = "E22"; Value = "  +6.12%  " }
    @{ Cell = "D23"; Value = "0.0000203" }
    @{ Cell = "E23"; Value = "  -3.33%  " }
    @{ Cell = "D24"; Value = "9.16" }
    @{ Cell = "E24"; Value = "  +5.70%  " }
    @{ Cell = "D25"; Value = "6.00" }
    @{ Cell = "E25"; Value = "  +8.24%  " }
    @{ Cell = "D26"; Value = "89.03" }
    @{ Cell = "E26"; Value = "  +3.04%  " }
    @{ Cell = "D27"; Value = "11.98" }
    @{ Cell = "E27"; Value = "  +3.10%  " }
    @{ Cell = "D29"; Value = "0.998" }
    @{ Cell = "E29"; Value = "  -0.24%  " }
    @{ Cell = "E30"; Value = "  +44.02%  " }
    @{ Cell = "D31"; Value = "0.171" }
    @{ Cell = "E31"; Value = "  +10.35%  " }
    @{ Cell = "E32"; Value = "  +17.83%  " }
    @{ Cell = "D33"; Value = "9.34" }
    @{ Cell = "E33"; Value = "  +4.30%  " }
    @{ Cell = "B34"; Value = "Binance-PegBSC-USD" }
    @{ Cell = "C34"; Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd" }
    @{ Cell = "D34"; Value = "0.999" }
    @{ Cell = "E34"; Value = "  -9.64%  " }
    @{ Cell = "B35"; Value = "Kaspa" }
    @{ Cell = "C35"; Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas" }
    @{ Cell = "D35"; Value = "0.170" }
    @{ Cell = "E35"; Value = "  +13.58%  " }
    @{ Cell = "B36"; Value = "EthereumClassic" }
    @{ Cell = "C36"; Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc" }
    @{ Cell = "D36"; Value = "26.47" }
    @{ Cell = "E36"; Value = "  +4.07%  " }
    @{ Cell = "B37"; Value = "RenderToken" }
    @{ Cell = "C37"; Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render" }
    @{ Cell = "D37"; Value = "7.62" }
    @{ Cell = "E37"; Value = "  +8.60%  " }
    @{ Cell = "D38"; Value = "511.66" }
    @{ Cell = "E38"; Value = "  +4.40%  " }
    @{ Cell = "D39"; Value = "1.96" }
    @{ Cell = "E40"; Value = "  +8.67%  " }
    @{ Cell = "D41"; Value = "0.451" }
    @{ Cell = "E41"; Value = "  +14.36%  " }
    @{ Cell = "D42"; Value = "3.83" }
    @{ Cell = "E42"; Value = "  +5.35%  " }
    @{ Cell = "D43"; Value = "3.44" }
    @{ Cell = "E43"; Value = "  -4.55%  " }
    @{ Cell = "D44"; Value = "22.17" }
    @{ Cell = "E44"; Value = "  +0.48%  " }
    @{ Cell = "B46"; Value = "Monero" }
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr" }
    @{ Cell = "D46"; Value = "157.95" }
    @{ Cell = "E46"; Value = "  +3.27%  " }
    @{ Cell = "B47"; Value = "ARBITRUM" }
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb" }
    @{ Cell = "D47"; Value = "0.710" }
    @{ Cell = "E47"; Value = "  +6.65%  " }
    @{ Cell = "E48"; Value = "  +5.54%  " }
    @{ Cell = "E49"; Value = "  +6.64%  " }
    @{ Cell = "D50"; Value = "4.48" }
    @{ Cell = "E50"; Value = "  +4.78%  " }
    @{ Cell = "D51"; Value = "44.00" }
    @{ Cell = "E51"; Value = "  -0.85%  " }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $value = $u.Value

    # Some "Price" values look like plain numbers (e.g. "238.66", "0.999",
    # "44.00"). Excel's normal type-inference would silently convert those
    # to numeric cells and drop meaningful trailing zeros / the text
    # formatting the source data relies on. Force the cell to text, write
    # the literal string, then drop the temporary number format so no
    # stray style is left behind on the cell.
    if ($value -match '^[+-]?\d+(\.\d+)?$') {
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.ClearFormats()
    } else {
        $range.Value = $value
    }
}
